$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "log_age_c" column header in G1 to "age_c" (value content change
# following feedback on the psychometrics paper).
$ws.Range("G1").Value = "age_c"

# Move the active selection to G1 (matches the saved selection state in the
# edited workbook).
$ws.Range("G1").Select()
